$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.694.11"
$ws.Range("E2").Value = "  +2.06%  "
$ws.Range("D3").Value = "2.937.30"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Formula = '="592.48"'
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("D6").Formula = '="146.69"'
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  +0.79%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "2.935.35"
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  +0.67%  "
$ws.Range("D10").Formula = '="7.29"'
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  +3.79%  "
$ws.Range("E11").Value = "  +5.15%  "
$ws.Range("D12").Formula = '="0.440"'
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("E13").Value = "  +4.45%  "
$ws.Range("D14").Formula = '="32.54"'
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  -3.07%  "
$ws.Range("E15").Value = "  -0.75%  "
$ws.Range("D16").Value = "3.420.71"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "62.646.58"
$ws.Range("E17").Value = "  +1.96%  "
$ws.Range("D18").Formula = '="6.66"'
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").Value = "2.949.57"
$ws.Range("E19").Value = "  +0.63%  "
$ws.Range("D20").Formula = '="438.46"'
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("D21").Formula = '="13.33"'
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  -1.06%  "
$ws.Range("D22").Formula = '="0.663"'
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  -1.96%  "
$ws.Range("D23").Formula = '="7.01"'
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  -1.23%  "
$ws.Range("D24").Formula = '="80.78"'
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  -1.44%  "
$ws.Range("D25").Formula = '="11.06"'
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  +1.64%  "
$ws.Range("D26").Formula = '="2.12"'
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  -3.55%  "
$ws.Range("D27").Formula = '="11.69"'
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("D28").Formula = '="1.00"'
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E29").Value = "  +0.56%  "
$ws.Range("E30").Value = "  +3.11%  "
$ws.Range("E31").Value = "  -0.41%  "
$ws.Range("D32").Formula = '="0.0000101"'
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  +13.74%  "
$ws.Range("E33").Value = "  -1.29%  "
$ws.Range("D34").Formula = '="26.26"'
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  -1.43%  "
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").Formula = '="0.989"'
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  -2.36%  "
$ws.Range("D37").Formula = '="3.08"'
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  +3.33%  "
$ws.Range("D38").Formula = '="5.55"'
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  -1.55%  "
$ws.Range("D39").Formula = '="49.64"'
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  -0.76%  "
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("D41").Formula = '="8.44"'
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  -1.68%  "
$ws.Range("D42").Formula = '="0.117"'
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  -5.03%  "
$ws.Range("D43").Formula = '="0.279"'
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  -1.31%  "
$ws.Range("D44").Formula = '="38.84"'
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  -8.30%  "
$ws.Range("D45").Value = "2.693.18"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("D46").Formula = '="134.09"'
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  -0.10%  "
$ws.Range("D47").Formula = '="360.47"'
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("D48").Formula = '="0.0335"'
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  -3.33%  "
$ws.Range("E50").Value = "  -1.05%  "
$ws.Range("D51").Formula = '="22.59"'
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  -4.58%  "
